$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rescatables")

$ids = @(20330051920306, 20330051920184, 20330051920136)
$paterno = @("PARRA", "VAZQUEZ", "MARCIAL")
$materno = @("FLORES", "VICTORIANO", "MORALES")
$nombres = @("SUEMI", "MARIAN", "IVAN DE JESUS")
$largo = "DISTINGUE LOS DIFERENTES TIPOS DE EMPRESA POR SU GIRO, ÁREAS FUNCIONALES, DOCUMENTACIÓN ADMINISTRATIVA Y RECURSOS"
$grupo = @("2ARHM", "2ARHM", "2ARHV")
$reprobadas = @(2, 2, 2)

for ($i = 0; $i -lt $ids.Count; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $ids[$i]
}
for ($i = 0; $i -lt $paterno.Count; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $paterno[$i]
}
for ($i = 0; $i -lt $materno.Count; $i++) {
    $ws.Cells.Item($i + 2, 3).Value = $materno[$i]
}
for ($i = 0; $i -lt $nombres.Count; $i++) {
    $ws.Cells.Item($i + 2, 4).Value = $nombres[$i]
}
for ($i = 0; $i -lt 3; $i++) {
    $ws.Cells.Item($i + 2, 5).Value = $largo
}
for ($i = 0; $i -lt $grupo.Count; $i++) {
    $ws.Cells.Item($i + 2, 6).Value = $grupo[$i]
}
for ($i = 0; $i -lt $reprobadas.Count; $i++) {
    $ws.Cells.Item($i + 2, 7).Value = $reprobadas[$i]
}
